# Update base file: add "Prior distribution" rows for fluxes / thermodynamic
# quantities right after the "LP solver" row, pushing the remaining
# parameter rows down by two, and move the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 6 ("Number of exp. conditions..."),
# shifting rows 6-12 down to rows 8-14.
$ws.Rows.Item(6).Resize(2).Insert()

# Copy the formatting (font/border/alignment) from the row above (row 5,
# "LP solver...") onto the two freshly inserted rows so they match the
# rest of the parameter table.
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8

# Populate the two new rows.
$ws.Range("A6").Value = "Prior distribution for fluxes (uniform or normal)"
$ws.Range("B6").Value = "normal"
$ws.Range("A7").Value = "Prior distribution for thermodynamic quantities (uniform or normal)"
$ws.Range("B7").Value = "normal"

# Move the active selection to A11:B12 (active cell B12), matching the
# updated worksheet view state.
$ws.Range("A11:B12").Select()
